{"js": "// Helper: find the single paragraph that contains `anchorText` (a text\n// fragment unique to that paragraph) and return the Word.Paragraph object.\nasync function findParagraph(context, anchorText) {\n  const results = context.document.body.search(anchorText, { matchCase: true });\n  results.load(\"paragraphs\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Paragraph anchor not found: \" + anchorText);\n  }\n  return results.items[0].paragraphs.items[0];\n}\n\n// --- 1. Remove the \"Jumping\" bullet entirely -------------------------------\nconst jumping = await findParagraph(context, \"Jumping\");\njumping.delete();\nawait context.sync();\n\n// --- 2. Remove the \"Collision detection ...\" bullet entirely ---------------\nconst collision = await findParagraph(context, \"Collision detection\");\ncollision.delete();\nawait context.sync();\n\n// --- 3. Textures bullet: drop \"background, \" and \"player, \" ----------------\nconst textures = await findParagraph(context, \"Textures (\");\ntextures.insertText(\n  \"Textures (enemies, objects, ground, etc.)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 4. Health system bullet: rewrite wording -------------------------------\nconst health = await findParagraph(context, \"Health system\");\nhealth.insertText(\n  \"Health system (if you have a power-up and get hit, you lose it. If you get hit without a power-up, you die)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 5. Remove the \"Kill enemies ...\" bullet entirely -----------------------\nconst killEnemies = await findParagraph(context, \"Kill enemies\");\nkillEnemies.delete();\nawait context.sync();\n\n// --- 6. Points system bullet: rewrite wording -------------------------------\nconst points = await findParagraph(context, \"Points system\");\npoints.insertText(\n  \"More stuff for Points system (item collection, finish level, lives bonus when you finish, etc.)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 7. Lives bullet: trim down --------------------------------------------\nconst lives = await findParagraph(context, \"Lives (find\");\nlives.insertText(\"Lives (find extra lives)\", \"Replace\");\nawait context.sync();\n\n// --- 8. Game over bullet: rewrite wording -----------------------------------\nconst gameOver = await findParagraph(context, \"Game over\");\ngameOver.insertText(\n  \"Game over (need to print a screen for x seconds or something, then reset)\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 9. Sound effects bullet: drop \"and background music\" ------------------\nconst sound = await findParagraph(context, \"Sound effects\");\nsound.insertText(\"Sound effects\", \"Replace\");\nawait context.sync();\n\n// --- 10. Append 4 blank paragraphs + 1 indented blank paragraph ------------\ncontext.document.body.paragraphs.load(\"items\");\nawait context.sync();\nlet tail = context.document.body.paragraphs.items[\n  context.document.body.paragraphs.items.length - 1\n];\n\nfor (let i = 0; i < 4; i++) {\n  const blank = tail.insertParagraph(\"\", \"After\");\n  blank.style = \"Normal\";\n  await context.sync();\n  tail = blank;\n}\n\nconst indented = tail.insertParagraph(\"\", \"After\");\nindented.style = \"Normal\";\nindented.firstLineIndent = 36; // 36pt = 720 twips\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---- helper: delete the whole paragraph containing $anchor -----------------\nfunction Remove-ParagraphByAnchor([string]$anchor) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    $rng.Expand(4) | Out-Null   # wdParagraph\n    $rng.Delete()\n}\n\n# ---- helper: replace the text of the paragraph containing $anchor ----------\n# (keeps the paragraph mark / pPr / numbering intact, only swaps the text)\nfunction Set-ParagraphTextByAnchor([string]$anchor, [string]$newText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    $rng.Expand(4) | Out-Null      # wdParagraph - grabs the paragraph incl. mark\n    $rng.MoveEnd(1, -1) | Out-Null # wdCharacter - shrink back off the mark\n    $rng.Text = $newText\n}\n\n# --- 1. Remove the \"Jumping\" bullet entirely --------------------------------\nRemove-ParagraphByAnchor \"Jumping\"\n\n# --- 2. Remove the \"Collision detection ...\" bullet entirely ---------------\nRemove-ParagraphByAnchor \"Collision detection\"\n\n# --- 3. Textures bullet: drop \"background, \" and \"player, \" ----------------\nSet-ParagraphTextByAnchor \"Textures (\" \"Textures (enemies, objects, ground, etc.)\"\n\n# --- 4. Health system bullet: rewrite wording -------------------------------\nSet-ParagraphTextByAnchor \"Health system\" \"Health system (if you have a power-up and get hit, you lose it. If you get hit without a power-up, you die)\"\n\n# --- 5. Remove the \"Kill enemies ...\" bullet entirely -----------------------\nRemove-ParagraphByAnchor \"Kill enemies\"\n\n# --- 6. Points system bullet: rewrite wording -------------------------------\nSet-ParagraphTextByAnchor \"Points system\" \"More stuff for Points system (item collection, finish level, lives bonus when you finish, etc.)\"\n\n# --- 7. Lives bullet: trim down ---------------------------------------------\nSet-ParagraphTextByAnchor \"Lives (find\" \"Lives (find extra lives)\"\n\n# --- 8. Game over bullet: rewrite wording ------------------------------------\nSet-ParagraphTextByAnchor \"Game over\" \"Game over (need to print a screen for x seconds or something, then reset)\"\n\n# --- 9. Sound effects bullet: drop \"and background music\" -------------------\nSet-ParagraphTextByAnchor \"Sound effects\" \"Sound effects\"\n\n# --- 10. Append 4 blank paragraphs + 1 indented blank paragraph ------------\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\nfor ($i = 0; $i -lt 4; $i++) {\n    $last.Range.InsertParagraphAfter()\n    $last = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $last.Style = \"Normal\"\n}\n$last.Range.InsertParagraphAfter()\n$indented = $d.Paragraphs.Item($d.Paragraphs.Count)\n$indented.Style = \"Normal\"\n$indented.Format.FirstLineIndent = 36\n"}
